# Applies odds-column corrections to rows 3, 7, 15, 16, 17, 18, 19
# per the commit diff (FlashScore odds workbook refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.55
$ws.Range("I3").Value = 5.75
$ws.Range("J3").Value = 2.1
$ws.Range("AQ3").Value = 23
$ws.Range("AX3").Value = 7

# Row 7
$ws.Range("H7").Value = 3.45
$ws.Range("I7").Value = 3.25
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 1.95
$ws.Range("U7").Value = 1.62
$ws.Range("V7").Value = 2.05
$ws.Range("W7").Value = 8.5
$ws.Range("X7").Value = 10.5
$ws.Range("AC7").Value = 11.25
$ws.Range("AF7").Value = 55
$ws.Range("AI7").Value = 18
$ws.Range("AP7").Value = 17.5
$ws.Range("AT7").Value = 2.8

# Row 15
$ws.Range("G15").Value = 2.18
$ws.Range("I15").Value = 3.2
$ws.Range("J15").Value = 2.8
$ws.Range("L15").Value = 3.65
$ws.Range("S15").Value = 1.38
$ws.Range("T15").Value = 2.47
$ws.Range("U15").Value = 1.84
$ws.Range("V15").Value = 1.86
$ws.Range("W15").Value = 5.9
$ws.Range("X15").Value = 8.5
$ws.Range("Y15").Value = 7.5
$ws.Range("Z15").Value = 17.5
$ws.Range("AA15").Value = 15.5
$ws.Range("AB15").Value = 23
$ws.Range("AC15").Value = 8.25
$ws.Range("AG15").Value = 300
$ws.Range("AH15").Value = 8.25
$ws.Range("AI15").Value = 14.5
$ws.Range("AL15").Value = 22
$ws.Range("AM15").Value = 26
$ws.Range("AN15").Value = 4.1
$ws.Range("AO15").Value = 11.75
$ws.Range("AP15").Value = 19.5
$ws.Range("AQ15").Value = 50
$ws.Range("AR15").Value = 80
$ws.Range("AS15").Value = 250
$ws.Range("AX15").Value = 5.2
$ws.Range("AY15").Value = 17.5
$ws.Range("AZ15").Value = 22
$ws.Range("BA15").Value = 80
$ws.Range("BB15").Value = 100
$ws.Range("BC15").Value = 250

# Row 16
$ws.Range("M16").Value = 1.01
$ws.Range("N16").Value = 16.5
$ws.Range("P16").Value = 7
$ws.Range("Q16").Value = 1.15
$ws.Range("R16").Value = 4.05
$ws.Range("U16").Value = 1.23
$ws.Range("V16").Value = 3.97

# Row 17
$ws.Range("I17").Value = 3.15
$ws.Range("J17").Value = 2.22
$ws.Range("K17").Value = 2.62
$ws.Range("L17").Value = 3.35
$ws.Range("AA17").Value = 10.75
$ws.Range("AB17").Value = 13
$ws.Range("AH17").Value = 18
$ws.Range("AN17").Value = 4.6
$ws.Range("AP17").Value = 11.75
$ws.Range("AT17").Value = 4.55
$ws.Range("AV17").Value = 29
$ws.Range("AX17").Value = 6.2
$ws.Range("AZ17").Value = 15.5
$ws.Range("BA17").Value = 55
$ws.Range("BB17").Value = 55
$ws.Range("BC17").Value = 100

# Row 18
$ws.Range("H18").Value = 5.4
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 5.3
$ws.Range("P18").Value = 8.1
$ws.Range("R18").Value = 3.72
$ws.Range("X18").Value = 9.75
$ws.Range("Z18").Value = 10
$ws.Range("AH18").Value = 30
$ws.Range("AN18").Value = 4.1
$ws.Range("AP18").Value = 10.25
$ws.Range("AQ18").Value = 12
$ws.Range("AR18").Value = 21
$ws.Range("AS18").Value = 70
$ws.Range("AT18").Value = 5.3
$ws.Range("AU18").Value = 7.1
$ws.Range("AV18").Value = 35
$ws.Range("AW18").Value = 350
$ws.Range("AX18").Value = 9.25
$ws.Range("AZ18").Value = 23
$ws.Range("BA18").Value = 150
$ws.Range("BB18").Value = 110
$ws.Range("BC18").Value = 150

# Row 19
$ws.Range("G19").Value = 1.08
$ws.Range("H19").Value = 8.25
$ws.Range("J19").Value = 1.35
$ws.Range("K19").Value = 3.15
$ws.Range("M19").Value = 1.02
$ws.Range("N19").Value = 10
$ws.Range("O19").Value = 1.12
$ws.Range("P19").Value = 5.3
$ws.Range("Q19").Value = 1.39
$ws.Range("R19").Value = 2.77
$ws.Range("S19").Value = 1.23
$ws.Range("T19").Value = 3.8
$ws.Range("U19").Value = 3
$ws.Range("V19").Value = 1.34
$ws.Range("W19").Value = 8.5
$ws.Range("X19").Value = 5.6
$ws.Range("Y19").Value = 14
$ws.Range("Z19").Value = 5.3
$ws.Range("AA19").Value = 13.5
$ws.Range("AB19").Value = 65
$ws.Range("AC19").Value = 10
$ws.Range("AD19").Value = 22
$ws.Range("AE19").Value = 65
$ws.Range("AF19").Value = 400
$ws.Range("AH19").Value = 90
$ws.Range("AI19").Value = 500
$ws.Range("AJ19").Value = 120
$ws.Range("AM19").Value = 500
$ws.Range("AN19").Value = 2.9
$ws.Range("AO19").Value = 3.95
$ws.Range("AP19").Value = 19.5
$ws.Range("AQ19").Value = 7.8
$ws.Range("AR19").Value = 40
$ws.Range("AS19").Value = 350
$ws.Range("AT19").Value = 3.8
$ws.Range("AU19").Value = 14
$ws.Range("AV19").Value = 175
$ws.Range("AX19").Value = 23
$ws.Range("AY19").Value = 250
$ws.Range("AZ19").Value = 150

Write-Output "Applied odds updates to rows 3, 7, 15, 16, 17, 18, 19"
